# Update three video lectures (slides 10, 11, 12) so the sample dictionary
# literal's insertion order matches the already-sorted-by-value output shown
# later in the code listing, and keep the derived tuple-list example on
# slide 12 consistent with that new order.

$p = $ppt.ActivePresentation

# --- Slide 10: "Sorting Lists of Tuples" ---
$shape10 = $p.Slides.Item(10).Shapes.Item(3)
$tr10 = $shape10.TextFrame.TextRange
$run10 = $tr10.Find(" = {'a':10, 'b':1, 'c':22}", 0)
$run10.Text = " = {'a':10, 'c':22, 'b':1}"

# --- Slide 11: "Using sorted()" ---
$shape11 = $p.Slides.Item(11).Shapes.Item(3)
$tr11 = $shape11.TextFrame.TextRange
$run11 = $tr11.Find(" = {'a':10, 'b':1, 'c':22}", 0)
$run11.Text = " = {'a':10 , 'c':22, 'b':1}"

# --- Slide 12: "Sort by Values Instead of Key" ---
$shape12 = $p.Slides.Item(12).Shapes.Item(3)
$tr12 = $shape12.TextFrame.TextRange
$run12dict = $tr12.Find(" = {'a':10, 'b':1, 'c':22}", 0)
$run12dict.Text = " = {'a':10, 'c':22, 'b':1}"

$run12tuple = $tr12.Find("[(10, 'a'), (1, 'b'), (22, 'c')]", 0)
$run12tuple.Text = "[(10, 'a') , (22, 'c'), (1, 'b')]"
